$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the x_fixStart, x_corrSteps, y_corrSteps, y_nrSteps and alienID
# values between trial row 2 and trial row 3 (columns B, D, E, G, H).
$ws.Range("B2").Value = 2
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 8
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 11

$ws.Range("B3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 16
